$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($row, $col, $value) {
    $t.Cell($row, $col).Range.Text = [string]$value
}

function Replace-InCell($row, $col, $oldText, $newText) {
    $containerRange = $t.Cell($row, $col).Range
    $full = $containerRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw ("text not found in cell (" + $row + "," + $col + "): " + $oldText)
    }
    $start = $containerRange.Start + $idx
    $end = $start + $oldText.Length
    $rr = $d.Range($start, $end)
    $rr.Text = $newText
}

# Row 5: Problem 4 / Part A -> Problem 4 / Part "-"
Set-CellText 5 2 "-"

# Row 6: Problem 4 / Part B -> Problem 5 / Part "-"
Set-CellText 6 1 "5"
Set-CellText 6 2 "-"

# Row 7: Problem 5 -> Problem 6 (Part already "-")
Set-CellText 7 1 "6"

# Row 8: Problem 6 -> Problem 7 (Part already "-")
Set-CellText 8 1 "7"

# Row 9: Problem 7 / Part A -> Problem 8 / Part "-"
Set-CellText 9 1 "8"
Set-CellText 9 2 "-"

# Row 10: Problem 7 / Part B -> Problem 9 / Part "-", plus solution text rewording
Set-CellText 10 1 "9"
Set-CellText 10 2 "-"

Replace-InCell 10 3 `
    "normality assumption is met. This means you need to do a qq-plot for each of the" `
    "normality assumption is met. This means you need to create a histogram for each"

Replace-InCell 10 3 `
    "groups. The qq-plots show that the groups are not perfectly normal, but they are" `
    "of the groups. The histograms show that the groups are not perfectly normal, but"

Replace-InCell 10 3 `
    "probably close enough to proceed with ANOVA." `
    "they are probably close enough to proceed with ANOVA."

# Row 11: Problem 7 / Part C -> Problem 10 / Part "-"
Set-CellText 11 1 "10"
Set-CellText 11 2 "-"

# Row 12: Problem 7 / Part D -> Problem 11 / Part "-"
Set-CellText 12 1 "11"
Set-CellText 12 2 "-"

# Row 13: Problem 7 / Part E -> Problem 12 / Part "-"
Set-CellText 13 1 "12"
Set-CellText 13 2 "-"

# Row 14: Problem 7 / Part F -> Problem 13 / Part "-"
Set-CellText 14 1 "13"
Set-CellText 14 2 "-"

# Row 15: Problem 7 / Part G -> Problem 14 / Part "-"
Set-CellText 15 1 "14"
Set-CellText 15 2 "-"

# Row 16: Problem 8 / Part (empty, unformatted paragraph) -> Problem 15 / Part "-"
Set-CellText 16 1 "15"

$partXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Compact"/><w:jc w:val="left"/></w:pPr><w:r><w:t xml:space="preserve">-</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$t.Cell(16, 2).Range.InsertXML($partXml)
